# correção das notas do fórum para matc65 em 2021.2
# Zera todas as notas/contagens (colunas B:J) que ainda possuem valores
# diferentes de zero, mantendo a coluna A (matricula) intacta.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Dados começam na linha 2 (linha 1 é o cabeçalho) e vão até a última
# linha usada; colunas B (2) até J (10).
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value
        if ($val -ne $null -and $val -ne 0) {
            $cell.Value = 0
        }
    }
}
